$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking values (e.g. "1.010", "30.431.82")
# are not auto-converted to numbers and lose their original formatting/precision.
$ws.Range("D2:D51").NumberFormat = "@"

# Apply updated cell values as described by the diff
$ws.Range('D2').Value = '30.431.82'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.107.90'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.95%  '
$ws.Range('D5').Value = '334.83'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').Value = '0.5218'
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('D8').Value = '0.4542'
$ws.Range('E8').Value = '  +5.34%  '
$ws.Range('D9').Value = '53.08'
$ws.Range('E9').Value = '  +16.04%  '
$ws.Range('D10').Value = '0.08922'
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('D11').Value = '1.173'
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('D12').Value = '24.27'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').Value = '2.108.12'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('D14').Value = '6.832'
$ws.Range('E14').Value = '  +2.97%  '
$ws.Range('D15').Value = '8.021'
$ws.Range('E15').Value = '  +4.96%  '
$ws.Range('D16').Value = '96.40'
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.00001140'
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('D19').Value = '0.06653'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = '19.18'
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('D22').Value = '6.350'
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('D23').Value = '30.525.88'
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').Value = '12.40'
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').Value = '2.367'
$ws.Range('E25').Value = '  +3.99%  '
$ws.Range('D26').Value = '2.352.00'
$ws.Range('E26').Value = '  +1.71%  '
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '163.46'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '2.536'
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').Value = '133.31'
$ws.Range('E30').Value = '  +2.14%  '
$ws.Range('D31').Value = '1.219'
$ws.Range('E31').Value = '  +2.74%  '
$ws.Range('D32').Value = '0.1071'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('D33').Value = '1.666'
$ws.Range('E33').Value = '  +6.46%  '
$ws.Range('D34').Value = '6.300'
$ws.Range('E34').Value = '  +3.85%  '
$ws.Range('D35').Value = '3.945'
$ws.Range('E35').Value = '  +3.02%  '
$ws.Range('D36').Value = '10.43'
$ws.Range('E36').Value = '  +8.88%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02583'
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.652'
$ws.Range('E38').Value = '  +4.84%  '
$ws.Range('D39').Value = '0.06828'
$ws.Range('E39').Value = '  +3.42%  '
$ws.Range('D40').Value = '0.2303'
$ws.Range('E40').Value = '  +3.21%  '
$ws.Range('D41').Value = '12.63'
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('D42').Value = '0.6863'
$ws.Range('E42').Value = '  +1.36%  '
$ws.Range('D43').Value = '1.251'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '1.006'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '2.334'
$ws.Range('E45').Value = '  +6.78%  '
$ws.Range('D46').Value = '14.04'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('D47').Value = '0.6366'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').Value = '3.672'
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').Value = '1.248'
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('D50').Value = '0.3407'
$ws.Range('E50').Value = '  +25.14%  '
$ws.Range('D51').Value = '83.22'
